$wb = $excel.ActiveWorkbook

# Hunk 0: ALC!row31
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 335.33334
$ws.Range("I31").Value = 335.33334
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1006.00002
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -776.0000200000001

# Hunk 1: ALC!row52
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 976.1429000000001
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 976.1429000000001
$ws.Range("K52").Value = 0
$ws.Range("L52").ClearContents()
$ws.Range("M52").Value = 2928.4287
$ws.Range("N52").Value = -3248.4287

# Hunk 2: ALC!row132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1874.1578
$ws.Range("I132").Value = 1164.0857
$ws.Range("J132").Value = 10158.333
$ws.Range("K132").Value = 3492.2571
$ws.Range("L132").Value = 30474.999
$ws.Range("M132").Value = -962.2571000000003
$ws.Range("N132").Value = -35534.999

# Hunk 3: ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3084.64
$ws.Range("I137").Value = 3519.2632
$ws.Range("J137").Value = 1708.3334
$ws.Range("K137").Value = 10557.7896
$ws.Range("L137").Value = 5125.0002
$ws.Range("M137").Value = -8007.7896
$ws.Range("N137").Value = -10225.0002

# Hunk 4: ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 8843.965
$ws.Range("I138").Value = 603.92
$ws.Range("J138").Value = 15281.5
$ws.Range("K138").Value = 1811.76
$ws.Range("L138").Value = 45844.5
$ws.Range("M138").Value = 3328.24
$ws.Range("N138").Value = -56124.5

# Hunk 5: ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21260.115
$ws.Range("I32").Value = 22835.666
$ws.Range("K32").Value = 22835.666
$ws.Range("M32").Value = -22548.666

# Hunk 6: ARM!row61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7438.615
$ws.Range("I61").Value = 5871.16
$ws.Range("J61").Value = 10237.643
$ws.Range("K61").Value = 5871.16
$ws.Range("L61").Value = 10237.643
$ws.Range("M61").Value = -5659.16
$ws.Range("N61").Value = -10661.643

# Hunk 7: ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2078.8157
$ws.Range("I74").Value = 1654.0385
$ws.Range("J74").Value = 2999.1667
$ws.Range("K74").Value = 1654.0385
$ws.Range("L74").Value = 2999.1667
$ws.Range("M74").Value = -780.0385000000001
$ws.Range("N74").Value = -4747.1667

# Hunk 8: ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2078.8157
$ws.Range("I77").Value = 1654.0385
$ws.Range("J77").Value = 2999.1667
$ws.Range("K77").Value = 8270.192500000001
$ws.Range("L77").Value = 14995.8335
$ws.Range("M77").Value = -3902.192500000001
$ws.Range("N77").Value = -23731.8335

# Hunk 9: ARM!row88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2392.5334
$ws.Range("I88").Value = 1494.2858
$ws.Range("K88").Value = 1494.2858
$ws.Range("M88").Value = -1088.2858

# Hunk 10: ARM!row91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2392.5334
$ws.Range("I91").Value = 1494.2858
$ws.Range("K91").Value = 1494.2858
$ws.Range("M91").Value = -90.28580000000011

# Hunk 11: ARM!row121
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 79800
$ws.Range("J121").Value = 79800
$ws.Range("L121").Value = 79800
$ws.Range("N121").Value = -83294

# Hunk 12: ARM!row136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 7438.615
$ws.Range("I136").Value = 5871.16
$ws.Range("J136").Value = 10237.643
$ws.Range("K136").Value = 17613.48
$ws.Range("L136").Value = 30712.929
$ws.Range("M136").Value = -15063.48
$ws.Range("N136").Value = -35812.929

# Hunk 13: BSM!row12
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 999.5
$ws.Range("I12").Value = 999.5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 999.5
$ws.Range("L12").ClearContents()
$ws.Range("M12").Value = -831.5
$ws.Range("N12").Value = 0

# Hunk 14: BSM!row22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 480.2857
$ws.Range("I22").Value = 393.33334
$ws.Range("J22").Value = 1002
$ws.Range("K22").Value = 393.33334
$ws.Range("L22").Value = 1002
$ws.Range("M22").Value = -220.33334
$ws.Range("N22").Value = -1348

# Hunk 15: BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2704.0322
$ws.Range("I134").Value = 2182.2354
$ws.Range("J134").Value = 3337.6428
$ws.Range("K134").Value = 6546.706200000001
$ws.Range("L134").Value = 10012.9284
$ws.Range("M134").Value = -4011.706200000001
$ws.Range("N134").Value = -15082.9284

# Hunk 16: CRP!row21
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 5554.5
$ws.Range("J21").Value = 5554.5
$ws.Range("L21").Value = 5554.5
$ws.Range("N21").Value = -6024.5

# Hunk 17: CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14729.037
$ws.Range("I31").Value = 14248.267
$ws.Range("J31").Value = 15330
$ws.Range("K31").Value = 14248.267
$ws.Range("L31").Value = 15330
$ws.Range("M31").Value = -13953.267
$ws.Range("N31").Value = -15920

# Hunk 18: CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 14729.037
$ws.Range("I34").Value = 14248.267
$ws.Range("J34").Value = 15330
$ws.Range("K34").Value = 14248.267
$ws.Range("L34").Value = 15330
$ws.Range("M34").Value = -14046.267
$ws.Range("N34").Value = -15734

# Hunk 19: CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2459043.2
$ws.Range("I58").Value = 4785981.5
$ws.Range("J58").Value = 2830.4443
$ws.Range("K58").Value = 4785981.5
$ws.Range("L58").Value = 2830.4443
$ws.Range("M58").Value = -4785778.5
$ws.Range("N58").Value = -3236.4443

# Hunk 20: CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2697.1482
$ws.Range("I132").Value = 2015.6
$ws.Range("J132").Value = 4644.4287
$ws.Range("K132").Value = 6046.799999999999
$ws.Range("L132").Value = 13933.2861
$ws.Range("M132").Value = -3516.799999999999
$ws.Range("N132").Value = -18993.2861

# Hunk 21: CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2459043.2
$ws.Range("I136").Value = 4785981.5
$ws.Range("J136").Value = 2830.4443
$ws.Range("K136").Value = 14357944.5
$ws.Range("L136").Value = 8491.332900000001
$ws.Range("M136").Value = -14355394.5
$ws.Range("N136").Value = -13591.3329

# Hunk 22: CUL!row5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7251743.5
$ws.Range("I5").Value = 744.6429000000001
$ws.Range("J5").Value = 18531076
$ws.Range("K5").Value = 2233.9287
$ws.Range("L5").Value = 55593228
$ws.Range("M5").Value = -2121.9287
$ws.Range("N5").Value = -55593452

# Hunk 23: CUL!row20
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 1547.1177
$ws.Range("I20").Value = 716.8333
$ws.Range("K20").Value = 2150.4999
$ws.Range("M20").Value = -1923.4999

# Hunk 24: CUL!row63
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3007.0715
$ws.Range("I63").Value = 1749.4
$ws.Range("J63").Value = 3705.7778
$ws.Range("K63").Value = 5248.200000000001
$ws.Range("L63").Value = 11117.3334
$ws.Range("M63").Value = -4499.200000000001
$ws.Range("N63").Value = -12615.3334

# Hunk 25: CUL!row66
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 3007.0715
$ws.Range("I66").Value = 1749.4
$ws.Range("J66").Value = 3705.7778
$ws.Range("K66").Value = 15744.6
$ws.Range("L66").Value = 33352.00019999999
$ws.Range("M66").Value = -12000.6
$ws.Range("N66").Value = -40840.00019999999

# Hunk 26: CUL!row120
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 18644
$ws.Range("I120").Value = 11333
$ws.Range("J120").Value = 25955
$ws.Range("K120").Value = 33999
$ws.Range("L120").Value = 77865
$ws.Range("M120").Value = -29161
$ws.Range("N120").Value = -87541

# Hunk 27: CUL!row122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1539.8636
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 1437.6111
$ws.Range("K122").Value = 18000
$ws.Range("L122").Value = 12938.4999
$ws.Range("M122").Value = -15550
$ws.Range("N122").Value = -17838.4999

# Hunk 28: CUL!row132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1932.2222
$ws.Range("I132").Value = 2379
$ws.Range("J132").Value = 1373.75
$ws.Range("K132").Value = 21411
$ws.Range("L132").Value = 12363.75
$ws.Range("M132").Value = -18881
$ws.Range("N132").Value = -17423.75

# Hunk 29: CUL!row135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 7251743.5
$ws.Range("I135").Value = 744.6429000000001
$ws.Range("J135").Value = 18531076
$ws.Range("K135").Value = 6701.7861
$ws.Range("L135").Value = 166779684
$ws.Range("M135").Value = -4166.7861
$ws.Range("N135").Value = -166784754

# Hunk 30: CUL!row139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 3161.4348
$ws.Range("I139").Value = 2152
$ws.Range("J139").Value = 3937.923
$ws.Range("K139").Value = 6456
$ws.Range("L139").Value = 11813.769
$ws.Range("M139").Value = -1316
$ws.Range("N139").Value = -22093.769

# Hunk 31: GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3385.2
$ws.Range("I132").Value = 2419.4167
$ws.Range("J132").Value = 4276.6924
$ws.Range("K132").Value = 7258.250100000001
$ws.Range("L132").Value = 12830.0772
$ws.Range("M132").Value = -4728.250100000001
$ws.Range("N132").Value = -17890.0772

# Hunk 32: LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4528
$ws.Range("I132").Value = 4182.4287
$ws.Range("J132").Value = 5737.5
$ws.Range("K132").Value = 12547.2861
$ws.Range("L132").Value = 17212.5
$ws.Range("M132").Value = -10017.2861
$ws.Range("N132").Value = -22272.5

# Hunk 33: LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4853.3687
$ws.Range("I136").Value = 2583.7778
$ws.Range("J136").Value = 6896
$ws.Range("K136").Value = 7751.3334
$ws.Range("L136").Value = 20688
$ws.Range("M136").Value = -5201.3334
$ws.Range("N136").Value = -25788

# Hunk 34: WVR!row32
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 19999.5
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 19999.5
$ws.Range("K32").Value = 0
$ws.Range("L32").ClearContents()
$ws.Range("M32").Value = 19999.5
$ws.Range("N32").Value = -20633.5

# Hunk 35: WVR!row96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3689
$ws.Range("J96").Value = 4429
$ws.Range("L96").Value = 4429
$ws.Range("N96").Value = -7175

# Hunk 36: WVR!row100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1528.8462
$ws.Range("I100").Value = 597.2
$ws.Range("J100").Value = 4634.3335
$ws.Range("K100").Value = 1194.4
$ws.Range("L100").Value = 9268.666999999999
$ws.Range("M100").Value = -653.4000000000001
$ws.Range("N100").Value = -10350.667

# Hunk 37: WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3434.322
$ws.Range("I136").Value = 2984.5386
$ws.Range("J136").Value = 3788.697
$ws.Range("K136").Value = 8953.6158
$ws.Range("L136").Value = 11366.091
$ws.Range("M136").Value = -6403.6158
$ws.Range("N136").Value = -16466.091

